$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 'Price' column stores values as plain text (e.g. thousands-grouped
# with dots, or with a fixed number of decimals), not as numbers. Whenever
# the refreshed price would otherwise be auto-recognized by Excel as a
# number, mark the cell as Text first so the literal string is preserved.

$ws.Range("D2").Value = '57.265.02'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '3.055.54'
$ws.Range("E3").Value = '  +1.65%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '517.00'
$ws.Range("E5").Value = '  +2.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.24'
$ws.Range("E6").Value = '  +1.41%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.442'
$ws.Range("E8").Value = '  +1.75%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.22'
$ws.Range("E9").Value = '  -4.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.110'
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.378'
$ws.Range("E11").Value = '  +3.76%  '
$ws.Range("D12").Value = '3.581.83'
$ws.Range("E12").Value = '  +1.57%  '
$ws.Range("E13").Value = '  -2.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.88'
$ws.Range("E14").Value = '  +2.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000166'
$ws.Range("E15").Value = '  +2.80%  '
$ws.Range("D16").Value = '57.154.51'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.16'
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").Value = '3.052.24'
$ws.Range("E18").Value = '  +1.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.48'
$ws.Range("E19").Value = '  +5.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.13'
$ws.Range("E20").Value = '  +2.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '333.04'
$ws.Range("E21").Value = '  +2.12%  '
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.506'
$ws.Range("E23").Value = '  +1.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.62'
$ws.Range("E24").Value = '  +1.73%  '
$ws.Range("D25").Value = '3.180.36'
$ws.Range("E25").Value = '  +1.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.39%  '
$ws.Range("E27").Value = '  -0.74%  '
$ws.Range("D28").Value = '0.0₃0899'
$ws.Range("E28").Value = '  -2.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.74'
$ws.Range("E29").Value = '  +0.28%  '
$ws.Range("E30").Value = '  -1.93%  '
$ws.Range("E31").Value = '  +0.98%  '
$ws.Range("E32").Value = '  +1.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.81'
$ws.Range("E33").Value = '  +1.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.72'
$ws.Range("E34").Value = '  -0.55%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '152.11'
$ws.Range("E35").Value = '  -1.00%  '
$ws.Range("E36").Value = '  +1.19%  '
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.17'
$ws.Range("E38").Value = '  +2.77%  '
$ws.Range("E39").Value = '  -0.65%  '
$ws.Range("D40").Value = '3.092.71'
$ws.Range("E40").Value = '  +1.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '37.18'
$ws.Range("E41").Value = '  -1.56%  '
$ws.Range("E42").Value = '  +1.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.666'
$ws.Range("E44").Value = '  +2.46%  '
$ws.Range("D45").Value = '2.209.57'
$ws.Range("E45").Value = '  -0.40%  '
$ws.Range("E46").Value = '  -0.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.967'
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.40'
$ws.Range("E49").Value = '  +4.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0242'
$ws.Range("E50").Value = '  +1.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0173'
$ws.Range("E51").Value = '  +9.91%  '
